# edit.ps1
# Applies two changes to the document:
#  1. In the "Income, age, ethnicity, ... population of zipcode,etc" paragraph,
#     remove the spell/grammar-check proofing marks that wrap "zipcode,etc" and
#     merge the trailing text into a single run of "... population of zipcode,etc".
#  2. After the "What factors occur in high covid zipcodes?" paragraph, add a
#     blank paragraph followed by a new "Census Pull:" section listing the data
#     that will be collected.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: clean up "Income, age, ethnicity, population of zipcode,etc"
# ---------------------------------------------------------------------------
# Locate the paragraph by its distinctive leading text.
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs($i)
    if ($candidate.Range.Text -like "Income, age, ethnicity,*zipcode,etc*") {
        $targetPara = $candidate
        break
    }
}

if ($targetPara -ne $null) {
    $paraStart = $targetPara.Range.Start
    $paraEnd = $targetPara.Range.End

    # Delete the whole paragraph (including its paragraph mark). This removes
    # the old runs as well as the w:proofErr spell/grammar-check bookmarks
    # that were anchored inside it.
    $d.Range($paraStart, $paraEnd).Delete()

    # Re-insert a clean paragraph with the corrected, merged text.
    $d.Range($paraStart, $paraStart).InsertBefore("Income, age, ethnicity, population of zipcode,etc`r")
}

# ---------------------------------------------------------------------------
# Change 2: add "Census Pull:" section after the questions list
# ---------------------------------------------------------------------------
$factorsPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs($i)
    if ($candidate.Range.Text -like "What factors occur in high covid zipcodes?*") {
        $factorsPara = $candidate
    }
}
if ($factorsPara -ne $null) {
    $endPos = $factorsPara.Range.End

    # Insert the new paragraphs in a single call. "Z" is a temporary
    # placeholder character standing in for what should become an empty
    # paragraph; inserting it avoids the end-of-document Range collapsing
    # that would otherwise merge the blank paragraph away.
    $d.Range($endPos, $endPos).InsertAfter("`rZ`rCensus Pull:`rIncome, Age, Sex, Race, and Zipcode")

    # Remove the temporary placeholder character, leaving a true empty
    # paragraph behind.
    $blankPara = $d.Paragraphs($d.Paragraphs.Count - 2)
    $placeholderStart = $blankPara.Range.Start
    $d.Range($placeholderStart, $placeholderStart + 1).Delete()
}
